$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.712.35"
$ws.Range("E2").Value = "  -0.36%  "

$ws.Range("D3").Value = "2.291.43"
$ws.Range("E3").Value = "  -0.19%  "

$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "119.90"
$ws.Range("E5").Value = "  +5.82%  "

$ws.Range("D6").Value = "267.90"
$ws.Range("E6").Value = "  -0.73%  "

$ws.Range("E7").Value = "  +1.96%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").Value = "48.41"
$ws.Range("E10").Value = "  +0.49%  "

$ws.Range("D11").Value = "0.0942"
$ws.Range("E11").Value = "  -0.36%  "

$ws.Range("D12").Value = "9.24"
$ws.Range("E12").Value = "  +2.30%  "

$ws.Range("E13").Value = "  +0.99%  "

$ws.Range("D14").Value = "15.62"
$ws.Range("E14").Value = "  -1.38%  "

$ws.Range("D15").Value = "0.897"
$ws.Range("E15").Value = "  +4.95%  "

$ws.Range("D16").Value = "2.634.80"
$ws.Range("E16").Value = "  -0.23%  "

$ws.Range("D17").Value = "2.292.86"
$ws.Range("E17").Value = "  +0.20%  "

$ws.Range("D18").Value = "43.759.14"
$ws.Range("E18").Value = "  -0.08%  "

$ws.Range("E19").Value = "  -0.19%  "

$ws.Range("D20").Value = "7.08"
$ws.Range("E20").Value = "  +3.17%  "

$ws.Range("D21").Value = "72.55"
$ws.Range("E21").Value = "  +0.52%  "

$ws.Range("D22").Value = "2.48"
$ws.Range("E22").Value = "  +1.57%  "

$ws.Range("D23").Value = "236.40"
$ws.Range("E23").Value = "  +1.40%  "

$ws.Range("D24").Value = "9.75"
$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("D25").Value = "2.90"
$ws.Range("E25").Value = "  -4.36%  "

$ws.Range("D26").Value = "1.02"
$ws.Range("E26").Value = "  +1.57%  "

$ws.Range("D27").Value = "11.88"
$ws.Range("E27").Value = "  +2.23%  "

$ws.Range("D28").Value = "42.97"
$ws.Range("E28").Value = "  +2.85%  "

$ws.Range("E29").Value = "  +0.29%  "

$ws.Range("D30").Value = "2.26"
$ws.Range("E30").Value = "  -0.51%  "

$ws.Range("D31").Value = "173.52"
$ws.Range("E31").Value = "  -1.08%  "

$ws.Range("D32").Value = "21.79"
$ws.Range("E32").Value = "  +1.20%  "

$ws.Range("D33").Value = "0.0916"
$ws.Range("E33").Value = "  -2.00%  "

$ws.Range("D34").Value = "5.80"
$ws.Range("E34").Value = "  +1.47%  "

$ws.Range("E35").Value = "  +2.04%  "

$ws.Range("D36").Value = "0.0384"
$ws.Range("E36").Value = "  +5.49%  "

$ws.Range("D37").Value = "4.76"
$ws.Range("E37").Value = "  +2.28%  "

$ws.Range("D38").Value = "3.96"
$ws.Range("E38").Value = "  +2.85%  "

$ws.Range("E39").Value = "  -0.65%  "

$ws.Range("D40").Value = "2.58"
$ws.Range("E40").Value = "  +8.43%  "

$ws.Range("D41").Value = "14.33"
$ws.Range("E41").Value = "  +3.91%  "

$ws.Range("D42").Value = "75.09"
$ws.Range("E42").Value = "  -0.24%  "

$ws.Range("D43").Value = "0.240"
$ws.Range("E43").Value = "  -2.31%  "

$ws.Range("D44").Value = "5.99"
$ws.Range("E44").Value = "  -5.59%  "

$ws.Range("E45").Value = "  -0.17%  "

$ws.Range("E46").Value = "  -1.14%  "

$ws.Range("E47").Value = "  +3.59%  "

$ws.Range("D48").Value = "8.63"
$ws.Range("E48").Value = "  -2.23%  "

$ws.Range("B49").Value = "ordi"
$ws.Range("C49").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D49").Value = "73.57"
$ws.Range("E49").Value = "  +36.62%  "

$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.101"
$ws.Range("E50").Value = "  +1.03%  "

$ws.Range("D51").Value = "102.15"
$ws.Range("E51").Value = "  +0.54%  "
